$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per commit (GitHub Actions scheduled refresh).

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '48.034.69'
$ws.Cells.Item(2, 5).Value = '  +0.18%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.500.23'
$ws.Cells.Item(3, 5).Value = '  -0.47%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '320.06'
$ws.Cells.Item(5, 5).Value = '  -1.10%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '107.34'
$ws.Cells.Item(6, 5).Value = '  -2.30%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.27%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.03%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -3.50%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '39.45'
$ws.Cells.Item(10, 5).Value = '  -3.78%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '20.15'
$ws.Cells.Item(11, 5).Value = '  +7.48%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.70%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.16%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -2.43%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '2.890.80'
$ws.Cells.Item(15, 5).Value = '  -0.54%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '2.494.34'
$ws.Cells.Item(16, 5).Value = '  -0.62%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -2.67%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '47.880.44'
$ws.Cells.Item(18, 5).Value = '  +0.06%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -3.33%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.71'
$ws.Cells.Item(20, 5).Value = '  +0.68%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.90%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -1.88%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '277.35'
$ws.Cells.Item(23, 5).Value = '  +11.75%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '71.46'
$ws.Cells.Item(24, 5).Value = '  +0.74%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -0.69%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.11%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '25.57'
$ws.Cells.Item(27, 5).Value = '  -1.65%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '9.72'
$ws.Cells.Item(28, 5).Value = '  -3.44%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.139'
$ws.Cells.Item(29, 5).Value = '  -0.96%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '35.03'
$ws.Cells.Item(30, 5).Value = '  -0.37%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '2.09'
$ws.Cells.Item(31, 5).Value = '  -7.92%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '49.44'
$ws.Cells.Item(32, 5).Value = '  -0.68%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '19.48'
$ws.Cells.Item(33, 5).Value = '  -3.48%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.13%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -1.73%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.48%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -2.07%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.86%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.88'
$ws.Cells.Item(39, 5).Value = '  -3.62%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -1.04%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '120.57'
$ws.Cells.Item(41, 5).Value = '  +0.76%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.34%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '21.16'
$ws.Cells.Item(43, 5).Value = '  -6.74%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.09%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.005.77'
$ws.Cells.Item(45, 5).Value = '  +0.02%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +1.79%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.89%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.32%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '8.97'
$ws.Cells.Item(49, 5).Value = '  -1.07%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '5.15'
$ws.Cells.Item(50, 5).Value = '  -1.40%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '80.10'
$ws.Cells.Item(51, 5).Value = '  +2.50%  '
